$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the quantity (G) values that were previously blank, driving the
# dependent formulas in column I (and the SUM total in I32) to recalculate.
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1

$excel.Calculate()
